$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.967.53'
$ws.Range("E2").Value = '  +3.10%  '
$ws.Range("D3").Value = '3.564.21'
$ws.Range("E3").Value = '  +2.25%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '616.05'
$ws.Range("E5").Value = '  +7.17%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.32'
$ws.Range("E6").Value = '  +1.22%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.618'
$ws.Range("E7").Value = '  +2.60%  '
$ws.Range("D8").Value = '3.560.06'
$ws.Range("E8").Value = '  +2.32%  '
$ws.Range("E9").Value = '  -0.01%  '
$ws.Range("E10").Value = '  +6.43%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.28'
$ws.Range("E11").Value = '  +12.62%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.588'
$ws.Range("E12").Value = '  +2.03%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '46.75'
$ws.Range("E13").Value = '  +1.50%  '
$ws.Range("E14").Value = '  +3.01%  '
$ws.Range("D15").Value = '4.132.62'
$ws.Range("E15").Value = '  +2.15%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.42'
$ws.Range("E16").Value = '  -0.33%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '621.47'
$ws.Range("E17").Value = '  +1.08%  '
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '70.967.50'
$ws.Range("E18").Value = '  +3.07%  '
$ws.Range("B19").Value = 'WrappedEther'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D19").Value = '3.561.46'
$ws.Range("E19").Value = '  +2.41%  '
$ws.Range("E20").Value = '  -0.87%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.39'
$ws.Range("E21").Value = '  +1.37%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.883'
$ws.Range("E22").Value = '  +1.02%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.49'
$ws.Range("E23").Value = '  -13.79%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '15.75'
$ws.Range("E24").Value = '  +0.42%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '96.97'
$ws.Range("E25").Value = '  +0.62%  '
$ws.Range("E26").Value = '  +2.20%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  +0.09%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.62'
$ws.Range("E28").Value = '  +0.91%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.58'
$ws.Range("E29").Value = '  +3.84%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.09'
$ws.Range("E30").Value = '  -1.19%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.54'
$ws.Range("E31").Value = '  +1.66%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.08'
$ws.Range("E32").Value = '  -1.61%  '
$ws.Range("E33").Value = '  +0.24%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.99'
$ws.Range("E34").Value = '  +1.48%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '576.76'
$ws.Range("E35").Value = '  -7.94%  '
$ws.Range("B36").Value = 'dogwifhat'
$ws.Range("C36").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.64'
$ws.Range("E36").Value = '  +7.26%  '
$ws.Range("B37").Value = 'Hedera'
$ws.Range("C37").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.101'
$ws.Range("E37").Value = '  -0.24%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '10.89'
$ws.Range("E38").Value = '  +2.36%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0477'
$ws.Range("E39").Value = '  +8.42%  '
$ws.Range("B40").Value = 'OKB'
$ws.Range("C40").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '57.73'
$ws.Range("E40").Value = '  +2.42%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.999'
$ws.Range("E41").Value = '  -0.17%  '
$ws.Range("E42").Value = '  +6.34%  '
$ws.Range("D43").Value = '3.373.44'
$ws.Range("E43").Value = '  +1.26%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.321'
$ws.Range("E44").Value = '  -0.50%  '
$ws.Range("E45").Value = '  +9.73%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '33.06'
$ws.Range("E46").Value = '  +2.22%  '
$ws.Range("D47").Value = '0.0₃0707'
$ws.Range("E47").Value = '  +3.80%  '
$ws.Range("E48").Value = '  +3.89%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.131'
$ws.Range("E49").Value = '  +1.80%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '133.81'
$ws.Range("E50").Value = '  +2.14%  '
$ws.Range("E51").Value = '  +1.42%  '
